$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "early end at" -> "early stop at" label in I1
$ws.Range("I1").Value = "early stop at"

# K7 previously held the long "early stopping, a learning rate scheduler, ..." note;
# it now holds a short status note, and a new note row (K8) is appended below it.
$ws.Range("K8").Value = "early stopping, warm up, slow down, random sample, random split"
$ws.Range("K7").Value = "some trials, not successful"

# Update the active selection to match the saved view state.
$ws.Range("M18").Select()
